$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.233.52"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.902.01"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.61"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5353"
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3812"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9008"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.65"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.338"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.001"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.83"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008635"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "27.269.31"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "1.125.61"
$ws.Range("E20").Value = "  -40.86%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.021"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.76"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.513"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.78"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.289"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.36"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.735"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.72"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.806"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.786"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09255"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8300"
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05053"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.218"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.999"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.327"
$ws.Range("E36").Value = "  -3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5724"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.076"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.285"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.565"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.27"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1520"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4939"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.13"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.635"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.46"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06134"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.95"
$ws.Range("E51").Value = "  -1.55%  "
